$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.067679594307382
$ws.Range("D2").Value = 1.066407713249684
$ws.Range("E2").Value = 1.071452302287768
$ws.Range("F2").Value = 1.079542135599914
$ws.Range("I2").Value = 1.046038474011035
$ws.Range("J2").Value = 1.072623263039358
$ws.Range("K2").Value = 1.069118343834023
$ws.Range("L2").Value = 1.074149437883078
$ws.Range("M2").Value = 1.082217914197943

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.069114138464885
$ws.Range("D3").Value = 1.067485898513489
$ws.Range("E3").Value = 1.072692109069217
$ws.Range("F3").Value = 1.080775892446985
$ws.Range("I3").Value = 1.046348868933507
$ws.Range("J3").Value = 1.073711991790584
$ws.Range("K3").Value = 1.070011460488644
$ws.Range("L3").Value = 1.075204755659584
$ws.Range("M3").Value = 1.083268750247478

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.070041802431421
$ws.Range("D4").Value = 1.068182835352751
$ws.Range("E4").Value = 1.073494016840618
$ws.Range("F4").Value = 1.081573801282444
$ws.Range("I4").Value = 1.046548065051609
$ws.Range("J4").Value = 1.074415421533425
$ws.Range("K4").Value = 1.070588047754985
$ws.Range("L4").Value = 1.075886719912772
$ws.Range("M4").Value = 1.083947726777279

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.070431657631806
$ws.Range("D5").Value = 1.06847565750326
$ws.Range("E5").Value = 1.073831062954426
$ws.Range("F5").Value = 1.081909145960604
$ws.Range("I5").Value = 1.046631413012683
$ws.Range("J5").Value = 1.074710895026525
$ws.Range("K5").Value = 1.070830131727059
$ws.Range("L5").Value = 1.076173205838834
$ws.Range("M5").Value = 1.084232935555082

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.070497108268382
$ws.Range("D6").Value = 1.068524813729913
$ws.Range("E6").Value = 1.073887650118832
$ws.Range("F6").Value = 1.081965446253595
$ws.Range("I6").Value = 1.046645384420545
$ws.Range("J6").Value = 1.074760491852503
$ws.Range("K6").Value = 1.070870760388988
$ws.Range("L6").Value = 1.076221295747097
$ws.Range("M6").Value = 1.084280809783682

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.070047012219359
$ws.Range("D7").Value = 1.068186748721872
$ws.Range("E7").Value = 1.073498520761701
$ws.Range("F7").Value = 1.081578282548709
$ws.Range("I7").Value = 1.046549180298446
$ws.Range("J7").Value = 1.074419370636739
$ws.Range("K7").Value = 1.070591283720959
$ws.Range("L7").Value = 1.075890548780578
$ws.Range("M7").Value = 1.083951538662843

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.06816452851177
$ws.Range("D8").Value = 1.066772241854073
$ws.Range("E8").Value = 1.071871371027985
$ws.Range("F8").Value = 1.079959177103674
$ws.Range("I8").Value = 1.046143715684217
$ws.Range("J8").Value = 1.072991423434311
$ws.Range("K8").Value = 1.069420451353884
$ws.Range("L8").Value = 1.07450627491398
$ws.Range("M8").Value = 1.082573254714964

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.064842687390795
$ws.Range("D9").Value = 1.064274062646976
$ws.Range("E9").Value = 1.069001456812277
$ws.Range("F9").Value = 1.07710279820803
$ws.Range("I9").Value = 1.045416553610883
$ws.Range("J9").Value = 1.070467011085584
$ws.Range("K9").Value = 1.067347083910478
$ws.Range("L9").Value = 1.072060008845984
$ws.Range("M9").Value = 1.08013687617183

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.062624699915222
$ws.Range("D10").Value = 1.062604657164431
$ws.Range("E10").Value = 1.067086188287851
$ws.Range("F10").Value = 1.075196130388799
$ws.Range("I10").Value = 1.04492319059676
$ws.Range("J10").Value = 1.068778366064113
$ws.Range("K10").Value = 1.065957817945277
$ws.Range("L10").Value = 1.070424276395137
$ws.Range("M10").Value = 1.078507292869471

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.061663408848762
$ws.Range("D11").Value = 1.061880811749488
$ws.Range("E11").Value = 1.066256338498263
$ws.Range("F11").Value = 1.07436990608512
$ws.Range("I11").Value = 1.044707507780898
$ws.Range("J11").Value = 1.06804576642748
$ws.Range("K11").Value = 1.06535454990176
$ws.Range("L11").Value = 1.069714785383795
$ws.Range("M11").Value = 1.077800360891305

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.061306203454834
$ws.Range("D12").Value = 1.061611792551954
$ws.Range("E12").Value = 1.065948012381001
$ws.Range("F12").Value = 1.074062911940763
$ws.Range("I12").Value = 1.04462708380217
$ws.Range("J12").Value = 1.0677734312294
$ws.Range("K12").Value = 1.065130209748009
$ws.Range("L12").Value = 1.069451063685283
$ws.Range("M12").Value = 1.077537574279896

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.061382831610345
$ws.Range("D13").Value = 1.061669504958972
$ws.Range("E13").Value = 1.066014153197382
$ws.Range("F13").Value = 1.074128767707264
$ws.Range("I13").Value = 1.044644349057355
$ws.Range("J13").Value = 1.067831857850194
$ws.Range("K13").Value = 1.065178343267756
$ws.Range("L13").Value = 1.069507641321918
$ws.Range("M13").Value = 1.077593952034612

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06163388500994
$ws.Range("D14").Value = 1.061858577633596
$ws.Range("E14").Value = 1.066230853880431
$ws.Range("F14").Value = 1.074344531867579
$ws.Range("I14").Value = 1.044700866236536
$ws.Range("J14").Value = 1.068023259525471
$ws.Range("K14").Value = 1.0653360111874
$ws.Range("L14").Value = 1.069692989864042
$ws.Range("M14").Value = 1.077778642980712

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.061788548693427
$ws.Range("D15").Value = 1.061975051470382
$ws.Range("E15").Value = 1.066364359212258
$ws.Range("F15").Value = 1.074477458233664
$ws.Range("I15").Value = 1.044735647248933
$ws.Range("J15").Value = 1.068141159796478
$ws.Range("K15").Value = 1.065433121110703
$ws.Range("L15").Value = 1.069807164574881
$ws.Range("M15").Value = 1.077892410480758

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.062688478445374
$ws.Range("D16").Value = 1.062652675435858
$ws.Range("E16").Value = 1.067141251269106
$ws.Range("F16").Value = 1.075250950690679
$ws.Range("I16").Value = 1.044937461383469
$ws.Range("J16").Value = 1.068826956345889
$ws.Range("K16").Value = 1.065997818633488
$ws.Range("L16").Value = 1.070471337213379
$ws.Range("M16").Value = 1.078554181680111

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.06325273888295
$ws.Range("D17").Value = 1.063077465879554
$ws.Range("E17").Value = 1.067628431556767
$ws.Range("F17").Value = 1.075735972173208
$ws.Range("I17").Value = 1.045063503408981
$ws.Range("J17").Value = 1.069256759402474
$ws.Range("K17").Value = 1.066351579332152
$ws.Range("L17").Value = 1.070887629158736
$ws.Range("M17").Value = 1.078968939611039

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.063581777172038
$ws.Range("D18").Value = 1.063325144449754
$ws.Range("E18").Value = 1.067912545105213
$ws.Range("F18").Value = 1.076018816997507
$ws.Range("I18").Value = 1.045136823556443
$ws.Range("J18").Value = 1.069507320941726
$ws.Range("K18").Value = 1.0665577575025
$ws.Range("L18").Value = 1.071130328894869
$ws.Range("M18").Value = 1.079210734580222

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.063693956432597
$ws.Range("D19").Value = 1.063409580509603
$ws.Range("E19").Value = 1.068009412041063
$ws.Range("F19").Value = 1.076115249710843
$ws.Range("I19").Value = 1.045161790280513
$ws.Range("J19").Value = 1.069592733119341
$ws.Range("K19").Value = 1.066628031129198
$ws.Range("L19").Value = 1.071213063614063
$ws.Range("M19").Value = 1.079293159097327

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.063192207941634
$ws.Range("D20").Value = 1.063031899685452
$ws.Range("E20").Value = 1.067576166972854
$ws.Range("F20").Value = 1.075683940193394
$ws.Range("I20").Value = 1.045050000785791
$ws.Range("J20").Value = 1.069210659627972
$ws.Range("K20").Value = 1.066313641191711
$ws.Range("L20").Value = 1.070842977031087
$ws.Range("M20").Value = 1.078924453102396

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.06155995992659
$ws.Range("D21").Value = 1.061802904632118
$ws.Range("E21").Value = 1.06616704325222
$ws.Range("F21").Value = 1.074280997426543
$ws.Range("I21").Value = 1.044684231897626
$ws.Range("J21").Value = 1.067966902471656
$ws.Range("K21").Value = 1.065289589104666
$ws.Range("L21").Value = 1.069638414478436
$ws.Range("M21").Value = 1.077724261666395

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.060532893239586
$ws.Range("D22").Value = 1.061029312654418
$ws.Range("E22").Value = 1.065280589309959
$ws.Range("F22").Value = 1.073398344745313
$ws.Range("I22").Value = 1.044452465562613
$ws.Range("J22").Value = 1.06718365665533
$ws.Range("K22").Value = 1.06464422372204
$ws.Range("L22").Value = 1.068879985169137
$ws.Range("M22").Value = 1.076968491307627

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.061077439181402
$ws.Range("D23").Value = 1.061439492341992
$ws.Range("E23").Value = 1.065750562321906
$ws.Range("F23").Value = 1.073866310645374
$ws.Range("I23").Value = 1.044575499659264
$ws.Range("J23").Value = 1.067598989546023
$ws.Range("K23").Value = 1.064986487655508
$ws.Range("L23").Value = 1.069282145828863
$ws.Range("M23").Value = 1.077369250838983

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.063219559531654
$ws.Range("D24").Value = 1.063052489379091
$ws.Range("E24").Value = 1.067599783244402
$ws.Range("F24").Value = 1.075707451390623
$ws.Range("I24").Value = 1.045056102651841
$ws.Range("J24").Value = 1.069231490549623
$ws.Range("K24").Value = 1.066330784312364
$ws.Range("L24").Value = 1.070863153764762
$ws.Range("M24").Value = 1.07894455503089

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.065702045544787
$ws.Range("D25").Value = 1.064920586572348
$ws.Range("E25").Value = 1.069743735003673
$ws.Range("F25").Value = 1.077841652257169
$ws.Range("I25").Value = 1.045606051622896
$ws.Range("J25").Value = 1.07112062281017
$ws.Range("K25").Value = 1.067884325135805
$ws.Range("L25").Value = 1.072693275409669
$ws.Range("M25").Value = 1.080767665029963
